# Regenerate all penyata to follow new data and format
#
# This script reproduces the content-level changes from the commit:
#  - relabels the "Kali Pertama/Kedua/Ketiga/Keempat" checkpoints to
#    "Semakan Kali ..." everywhere they appear
#  - title-cases the competition entry names
#  - adds two new competition entries (Sarung Race / Theme Party) with
#    their merit values
#  - updates the "Semakan Kali Ketiga" merit/demerit figures
#  - adjusts merged ranges / label placement to match the new layout
#  - tweaks print setup (centered horizontally, fit to one page tall)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header area: move "STATEMENT OF HOMEROOM ACCOUNT" from E4 to D4
# and widen its merge so it spans D4:G4 -------------------------------
$ws.Range("E4").Value = ""
$ws.Range("D4").Value = "STATEMENT OF HOMEROOM ACCOUNT"
$ws.Range("D4:G4").Merge()

# Widen / add a few merges used by the new layout
$ws.Range("B5:C5").Merge()
$ws.Range("B12:F12").Merge()
$ws.Range("B15:C15").UnMerge()
$ws.Range("B21:E21").Merge()
$ws.Range("B27:E27").Merge()
$ws.Range("B33:E33").Merge()
$ws.Range("B43:E43").Merge()

# --- Rename "Kali Pertama/Kedua/Ketiga/Keempat" labels --------------
# Merit Pendahuluan filing section (rows 16-19)
$ws.Range("C16").Value = "Semakan Kali Pertama"
$ws.Range("C17").Value = "Semakan Kali Kedua"
$ws.Range("C18").Value = "Semakan Kali Ketiga"
$ws.Range("C19").Value = "Semakan Kali Keempat"

# Laporan Atas Talian section (rows 22-25)
$ws.Range("C22").Value = "Semakan Kali Pertama"
$ws.Range("C23").Value = "Semakan Kali Kedua"
$ws.Range("C24").Value = "Semakan Kali Ketiga"
$ws.Range("C25").Value = "Semakan Kali Keempat"

# JPPM / JDM / JDRM section (rows 28-31)
$ws.Range("C28").Value = "Semakan Kali Pertama"
$ws.Range("C29").Value = "Semakan Kali Kedua"
$ws.Range("C30").Value = "Semakan Kali Ketiga"
$ws.Range("C31").Value = "Semakan Kali Keempat"

# --- Updated merit figures for "Semakan Kali Ketiga" (row 18) -------
$ws.Range("D18").Value = 10484
$ws.Range("E18").Value = 700

# --- Penyertaan Pertandingan: title-case existing entries + two new -
$ws.Range("C34").Value = "Unga"
$ws.Range("C35").Value = "Bouquet Kreatif"
$ws.Range("C36").Value = "Tik Tok Raya"
$ws.Range("C37").Value = "Riang Ria Kuih Raya"
$ws.Range("C38").Value = "Creative Collage"

$ws.Range("C39").Value = "Sarung Race"
$ws.Range("D39").Value = 100

$ws.Range("C40").Value = "Theme Party"
$ws.Range("D40").Value = 100

# --- Print setup: center horizontally, fit to one page tall ---------
$ws.PageSetup.CenterHorizontally = $true
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
$ws.PageSetup.HeaderMargin = 0
$ws.PageSetup.FooterMargin = 0
